# Weekly price-sheet update: a new record is inserted at the top of the
# data block (row 112), pushing the existing rows 112-182 down to 113-183.
#
# Insert a new row above row 112. Excel shifts every row at/after 112 down
# by one (112->113, ..., 182->183), carrying over all values/styles, and
# the worksheet dimension grows from A1:T182 to A1:T183 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly record.
$ws.Range("A112").Value = 7
$ws.Range("B112").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C112").Value = "Ñuble"
$ws.Range("D112").Value = 45161
$ws.Range("D112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E112").Value = 16
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100108
$ws.Range("H112").Value = "Tropicales y subtropicales"
$ws.Range("I112").Value = 100108002
$ws.Range("J112").Value = "Mango"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 60
$ws.Range("N112").Value = 9000
$ws.Range("O112").Value = 9000
$ws.Range("P112").Value = 9000
$ws.Range("Q112").Value = '$/bandeja 4 kilos'
$ws.Range("R112").Value = "Brasil"
$ws.Range("S112").Value = 2250
$ws.Range("T112").Value = 4
